$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 17:08"

# Update country rows: new COVID-19 figures refresh changed the sort order
# (tied/updated totals), so several rows now show a different country name
# together with its own updated Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes figures.

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2507369
$ws.Range("C4").Value = 2781
$ws.Range("D4").Value = 1052449
$ws.Range("E4").Value = 1328074
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 66
$ws.Range("H4").Value = 126846

# Row 7: India
$ws.Range("A7").Value = "India"
$ws.Range("B7").Value = 497824
$ws.Range("C7").Value = 6654
$ws.Range("D7").Value = 289768
$ws.Range("E7").Value = 192650
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 98
$ws.Range("H7").Value = 15406

# Row 9: España
$ws.Range("A9").Value = "España"
$ws.Range("B9").Value = 294985
$ws.Range("C9").Value = 419
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 8
$ws.Range("H9").Value = 28338

# Row 31: Argentina
$ws.Range("A31").Value = "Argentina"
$ws.Range("B31").Value = 52457
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 18416
$ws.Range("E31").Value = 32874
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 17
$ws.Range("H31").Value = 1167

# Row 36: Singapur
$ws.Range("A36").Value = "Singapur"
$ws.Range("B36").Value = 42955
$ws.Range("C36").Value = 219
$ws.Range("D36").Value = 36825
$ws.Range("E36").Value = 6104
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 26

# Row 45: Republica Dominicana
$ws.Range("A45").Value = "Republica Dominicana"
$ws.Range("B45").Value = 29764
$ws.Range("C45").Value = 623
$ws.Range("D45").Value = 16357
$ws.Range("E45").Value = 12695
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 14
$ws.Range("H45").Value = 712

# Row 57: Moldavia
$ws.Range("A57").Value = "Moldavia"
$ws.Range("B57").Value = 15776
$ws.Range("C57").Value = 323
$ws.Range("D57").Value = 8765
$ws.Range("E57").Value = 6500
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 9
$ws.Range("H57").Value = 511

# Row 58: Guatemala
$ws.Range("A58").Value = "Guatemala"
$ws.Range("B58").Value = 15619
$ws.Range("C58").Value = 800
$ws.Range("D58").Value = 2949
$ws.Range("E58").Value = 12047
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 22
$ws.Range("H58").Value = 623

# Row 59: Ghana
$ws.Range("A59").Value = "Ghana"
$ws.Range("B59").Value = 15473
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 11431
$ws.Range("E59").Value = 3947
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 95

# Row 60: Azerbaiyan
$ws.Range("A60").Value = "Azerbaiyan"
$ws.Range("B60").Value = 15369
$ws.Range("C60").Value = 517
$ws.Range("D60").Value = 8364
$ws.Range("E60").Value = 6818
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 187

# Row 61: Honduras
$ws.Range("A61").Value = "Honduras"
$ws.Range("B61").Value = 15366
$ws.Range("C61").Value = 795
$ws.Range("D61").Value = 1600
$ws.Range("E61").Value = 13340
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 9
$ws.Range("H61").Value = 426

# Row 71: Noruega
$ws.Range("A71").Value = "Noruega"
$ws.Range("B71").Value = 8828
$ws.Range("C71").Value = 40
$ws.Range("D71").Value = 8138
$ws.Range("E71").Value = 441
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 249

# Row 82: Kenia
$ws.Range("A82").Value = "Kenia"
$ws.Range("B82").Value = 5533
$ws.Range("C82").Value = 149
$ws.Range("D82").Value = 1905
$ws.Range("E82").Value = 3491
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 5
$ws.Range("H82").Value = 137

# Row 83: El Salvador
$ws.Range("A83").Value = "El Salvador"
$ws.Range("B83").Value = 5517
$ws.Range("C83").Value = 181
$ws.Range("D83").Value = 3291
$ws.Range("E83").Value = 2093
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 7
$ws.Range("H83").Value = 133

# Row 84: Etiopia
$ws.Range("A84").Value = "Etiopia"
$ws.Range("B84").Value = 5425
$ws.Range("C84").Value = 250
$ws.Range("D84").Value = 1688
$ws.Range("E84").Value = 3648
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 8
$ws.Range("H84").Value = 89

# Row 103: Cuba
$ws.Range("A103").Value = "Cuba"
$ws.Range("B103").Value = 2325
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 2180
$ws.Range("E103").Value = 60
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 85

# Row 125: Tunez
$ws.Range("A125").Value = "Tunez"
$ws.Range("B125").Value = 1164
$ws.Range("C125").Value = 2
$ws.Range("D125").Value = 1023
$ws.Range("E125").Value = 91
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 50

# Row 132: Cabo Verde
$ws.Range("A132").Value = "Cabo Verde"
$ws.Range("B132").Value = 1027
$ws.Range("C132").Value = 24
$ws.Range("D132").Value = 562
$ws.Range("E132").Value = 456
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 9

# Row 142: Mozambique
$ws.Range("A142").Value = "Mozambique"
$ws.Range("B142").Value = 816
$ws.Range("C142").Value = 28
$ws.Range("D142").Value = 223
$ws.Range("E142").Value = 588
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 5

# Row 148: Liberia
$ws.Range("A148").Value = "Liberia"
$ws.Range("B148").Value = 684
$ws.Range("C148").Value = 3
$ws.Range("D148").Value = 285
$ws.Range("E148").Value = 365
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 34

# Row 149: Jamaica
$ws.Range("A149").Value = "Jamaica"
$ws.Range("B149").Value = 684
$ws.Range("C149").Value = 6
$ws.Range("D149").Value = 526
$ws.Range("E149").Value = 148
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 10

# Row 200: Laos
$ws.Range("A200").Value = "Laos"
$ws.Range("B200").Value = 19
$ws.Range("C200").Value = 0
$ws.Range("D200").Value = 19
$ws.Range("E200").Value = 0
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 0

# Row 201: Santa Lucia
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 19
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 19
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# Row 202: Dominica
$ws.Range("A202").Value = "Dominica"
$ws.Range("B202").Value = 18
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 18
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

# Row 203: Fiyi
$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 208: Islas Malvinas
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("B208").Value = 13
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 13
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209: Groenlandia
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 212: Montserrat
$ws.Range("A212").Value = "Montserrat"
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 10
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 1

# Row 213: Seychelles
$ws.Range("A213").Value = "Seychelles"
$ws.Range("B213").Value = 11
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 11
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0
